$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Investor 1"
$ws.Range("A3").Value = "Investor 2"

$ws.Range("A4").Select()
